$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9964644908905029
$ws.Range("B1").Value = 1.711386799812317
$ws.Range("C1").Value = 3.468760967254639
$ws.Range("D1").Value = 3.706749677658081
$ws.Range("E1").Value = 0.9452120065689087
